# Generate Report for Handback
#
# Updates the localization-status workbook to reflect a failed handback
# transform:
#   - Every "Status" cell that said "Ready for handoff" (the Overview
#     rollup plus each language sheet) now says "Handback transform
#     failed".
#   - Each language sheet gets an "Error Detail" (column L) entry on the
#     74fe6409-...-handoff row explaining the handback/handoff file name
#     mismatch.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhSheet = $wb.Worksheets.Item("zh-cn")
$deSheet = $wb.Worksheets.Item("de-de")

# Status text update - every cell that shared the old "Ready for
# handoff" string gets the new text so the shared string is effectively
# edited in place rather than duplicated.
$overview.Range("B3").Value = "Handback transform failed"
$overview.Range("C3").Value = "Handback transform failed"
$zhSheet.Range("C3").Value = "Handback transform failed"
$deSheet.Range("C3").Value = "Handback transform failed"

# Error Detail (column L) for row 3 on each language sheet.
$zhSheet.Range("L3").Value = "Handback file name: 2uqgpt5d.sk0 is different with handoff file name: 74fe6409-53e5-4a90-b69a-231ff2bcc816.d549c47b6b7fcdb552633ec25b4799f2aae8884b.zh-cn."
$deSheet.Range("L3").Value = "Handback file name: 2uqgpt5d.sk0 is different with handoff file name: 74fe6409-53e5-4a90-b69a-231ff2bcc816.d549c47b6b7fcdb552633ec25b4799f2aae8884b.de-de."
